$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "41.536.28", "4.30") that must
# stay literal text rather than become a Number. A leading apostrophe forces
# Excel's text-entry interpretation of .Value; Style is then reset to
# "Normal" so no stray number-format/quote-prefix style is left behind.

$ws.Range("D2").Value = '''41.536.28'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.31%  '
$ws.Range("D3").Value = '''2.163.59'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.94%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''238.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.13%  '
$ws.Range("D6").Value = '''0.606'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.31%  '
$ws.Range("D7").Value = '''71.92'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.72%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '''0.575'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.54%  '
$ws.Range("D10").Value = '''39.89'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.38%  '
$ws.Range("D11").Value = '''0.0905'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.63%  '
$ws.Range("D12").Value = '''54.04'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.01%  '
$ws.Range("E13").Value = '  -3.15%  '
$ws.Range("D14").Value = '''6.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.25%  '
$ws.Range("D15").Value = '''2.485.48'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = '''14.14'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.34%  '
$ws.Range("D17").Value = '''2.153.81'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.06%  '
$ws.Range("D18").Value = '''0.778'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -7.27%  '
$ws.Range("D19").Value = '''41.377.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.34%  '
$ws.Range("D20").Value = '''0.0000103'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.26%  '
$ws.Range("D21").Value = '''69.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.42%  '
$ws.Range("E22").Value = '  -7.47%  '
$ws.Range("D23").Value = '''9.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -13.52%  '
$ws.Range("D24").Value = '''226.99'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.67%  '
$ws.Range("D25").Value = '''2.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.65%  '
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("D27").Value = '''10.65'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.88%  '
$ws.Range("D28").Value = '''3.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.76%  '
$ws.Range("E29").Value = '  -4.50%  '
$ws.Range("E30").Value = '  -1.41%  '
$ws.Range("D31").Value = '''170.13'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.81%  '
$ws.Range("D32").Value = '''19.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.03%  '
$ws.Range("D33").Value = '''33.22'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +10.14%  '
$ws.Range("D34").Value = '''0.0768'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.09%  '
$ws.Range("E35").Value = '  -9.77%  '
$ws.Range("E36").Value = '  -3.85%  '
$ws.Range("D37").Value = '''4.30'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("E38").Value = '  -4.81%  '
$ws.Range("D39").Value = '''0.0302'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.57%  '
$ws.Range("D40").Value = '''2.09'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.94%  '
$ws.Range("D41").Value = '''11.93'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -11.04%  '
$ws.Range("D42").Value = '''5.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.46%  '
$ws.Range("D43").Value = '''58.96'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -10.00%  '
$ws.Range("E44").Value = '  -5.26%  '
$ws.Range("D45").Value = '''8.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.27%  '
$ws.Range("D46").Value = '''0.0960'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.00%  '
$ws.Range("D47").Value = '''96.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.03%  '
$ws.Range("E48").Value = '  -3.99%  '
$ws.Range("E49").Value = '  -5.54%  '
$ws.Range("D50").Value = '''2.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.74%  '
$ws.Range("E51").Value = '  -2.29%  '
